$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "nama"
$ws.Range("C1").Value = "level_id"

$ws.Range("P15").Select() | Out-Null
